$wb = $excel.ActiveWorkbook

$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("E4").Value = "2016-03-23 20:50:26"
$wsZh.Range("H4").Value = "2016-03-23 20:51:02"

$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("E4").Value = "2016-03-23 20:50:31"
$wsDe.Range("H4").Value = "2016-03-23 20:51:10"
